$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs, Gnai2, Cxcr2, ECs)
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 2.537827835546889
$ws.Range("R2").Value = 22.840450519922
$ws.Range("S2").Value = 0.4067926910433548
$ws.Range("T2").Value = 0.4067926910433549

# Row 3 (FAPs, Gnai2, Cxcr2, ECs)
$ws.Range("I3").Value = 0.3656254573230189
$ws.Range("J3").Value = 0.365625457323019
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03970866666666666
$ws.Range("N3").Value = 0.119126
$ws.Range("Q3").Value = 2.2810008228
$ws.Range("R3").Value = 20.5290074052
$ws.Range("S3").Value = 0.3656254573230189
$ws.Range("T3").Value = 0.365625457323019

# Row 4 (MuSCs, Gnai2, Cxcr2, ECs)
$ws.Range("G4").Value = 35.755375
$ws.Range("H4").Value = 107.266125
$ws.Range("I4").Value = 0.2275818516336261
$ws.Range("J4").Value = 0.2275818516336262
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03970866666666666
$ws.Range("N4").Value = 0.119126
$ws.Range("Q4").Value = 1.419798267416666
$ws.Range("R4").Value = 12.77818440675
$ws.Range("S4").Value = 0.2275818516336261
$ws.Range("T4").Value = 0.2275818516336262
